$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (week number + date range) ---
$ws.Range("A8").Value = "Volume 32   Number  22"
$ws.Range("C9").Value = "Report Covering the Week  5/26/2025  Through  6/1/2025"

# --- Data table updates ---
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50
$ws.Range("N15").Value = -9.090909090909
$ws.Range("D16").Value = 5
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -66.666666666666
$ws.Range("J16").Value = 49
$ws.Range("K16").Value = -14.285714285714
$ws.Range("M16").Value = 16.666666666666
$ws.Range("N16").Value = -83.2
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 12
$ws.Range("H17").Value = -14.285714285714
$ws.Range("I17").Value = 60
$ws.Range("J17").Value = 53
$ws.Range("K17").Value = 13.207547169811
$ws.Range("M17").Value = 1.694915254237
$ws.Range("N17").Value = -38.144329896907
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 43
$ws.Range("J18").Value = 36
$ws.Range("K18").Value = 19.444444444444
$ws.Range("L18").Value = -15.686274509803
$ws.Range("M18").Value = -6.521739130434
$ws.Range("N18").Value = -79.024390243902
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 40
$ws.Range("F19").Value = 46
$ws.Range("G19").Value = 52
$ws.Range("H19").Value = -11.538461538461
$ws.Range("I19").Value = 249
$ws.Range("J19").Value = 283
$ws.Range("K19").Value = -12.014134275618
$ws.Range("L19").Value = -17
$ws.Range("M19").Value = 1.632653061224
$ws.Range("N19").Value = -29.261363636363
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 9
$ws.Range("K20").Value = -57.142857142857
$ws.Range("L20").Value = -60.869565217391
$ws.Range("M20").Value = -50
$ws.Range("N20").Value = -95.964125560538
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 5
$ws.Range("G21").Value = 87
$ws.Range("H21").Value = -19.540229885057
$ws.Range("I21").Value = 413
$ws.Range("J21").Value = 448
$ws.Range("K21").Value = -7.8125
$ws.Range("L21").Value = -15.885947046843
$ws.Range("M21").Value = 1.225490196078
$ws.Range("N21").Value = -63.803680981595
$ws.Range("L22").Value = -29.411764705882
$ws.Range("D23").Value = 2
$ws.Range("G23").Value = 9
$ws.Range("J23").Value = 22
$ws.Range("K23").Value = -27.272727272727
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 16.666666666666
$ws.Range("F24").Value = 68
$ws.Range("G24").Value = 72
$ws.Range("H24").Value = -5.555555555555
$ws.Range("I24").Value = 351
$ws.Range("J24").Value = 275
$ws.Range("K24").Value = 27.636363636363
$ws.Range("L24").Value = 20.61855670103
$ws.Range("M24").Value = 2.03488372093
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 6
$ws.Range("F25").Value = 33
$ws.Range("H25").Value = 43.478260869565
$ws.Range("I25").Value = 192
$ws.Range("J25").Value = 114
$ws.Range("K25").Value = 68.421052631578
$ws.Range("L25").Value = 15.662650602409
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = -30
$ws.Range("F26").Value = 29
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 16
$ws.Range("I26").Value = 121
$ws.Range("J26").Value = 132
$ws.Range("K26").Value = -8.333333333333
$ws.Range("L26").Value = -16.551724137931
$ws.Range("M26").Value = -19.867549668874
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -80
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 50
$ws.Range("I28").Value = 20
$ws.Range("K28").Value = -13.043478260869
$ws.Range("L28").Value = -31.03448275862
$ws.Range("L31").Value = -84.615384615384

# --- Cells converted from numeric to placeholder text (style 13) ---
# Force text storage via "@" format, then paste the General-format/style
# from an existing placeholder cell (C14, style 13) so the resulting style
# matches exactly while the stored value remains literal text.
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("G22").PasteSpecial(-4122)

$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("H22").PasteSpecial(-4122)

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)

